$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking price values (column D)
# so they remain strings like the original inline-string cells, not numbers.
$textCells = @("D5", "D6", "D10", "D14", "D15", "D19", "D21", "D24", "D25", "D29", "D30", "D31", "D32", "D33", "D37", "D38", "D39", "D41", "D42", "D43", "D46", "D47", "D48", "D49", "D50")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = '42.689.30'
$ws.Range("E2").Value = '  +1.91%  '
$ws.Range("D3").Value = '2.295.52'
$ws.Range("E3").Value = '  +0.28%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '322.55'
$ws.Range("E5").Value = '  +1.74%  '
$ws.Range("D6").Value = '104.77'
$ws.Range("E6").Value = '  +1.90%  '
$ws.Range("E7").Value = '  +0.51%  '
$ws.Range("E8").Value = '  +0.19%  '
$ws.Range("E9").Value = '  +0.73%  '
$ws.Range("D10").Value = '40.29'
$ws.Range("E10").Value = '  +3.36%  '
$ws.Range("E11").Value = '  +0.14%  '
$ws.Range("E12").Value = '  +2.80%  '
$ws.Range("E13").Value = '  +0.38%  '
$ws.Range("D14").Value = '0.972'
$ws.Range("E14").Value = '  +0.61%  '
$ws.Range("D15").Value = '15.27'
$ws.Range("E15").Value = '  -0.12%  '
$ws.Range("D16").Value = '2.644.44'
$ws.Range("E16").Value = '  +0.31%  '
$ws.Range("D17").Value = '2.294.23'
$ws.Range("E17").Value = '  -0.29%  '
$ws.Range("D18").Value = '42.625.63'
$ws.Range("D19").Value = '7.47'
$ws.Range("E19").Value = '  -1.38%  '
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("D21").Value = '13.26'
$ws.Range("E21").Value = '  +33.46%  '
$ws.Range("E22").Value = '  -0.52%  '
$ws.Range("E23").Value = '  +0.25%  '
$ws.Range("D24").Value = '270.09'
$ws.Range("E24").Value = '  -5.17%  '
$ws.Range("D25").Value = '2.23'
$ws.Range("E25").Value = '  -1.41%  '
$ws.Range("E26").Value = '  -0.34%  '
$ws.Range("E27").Value = '  +1.38%  '
$ws.Range("E28").Value = '  +3.09%  '
$ws.Range("D29").Value = '38.67'
$ws.Range("E29").Value = '  +10.88%  '
$ws.Range("D30").Value = '22.61'
$ws.Range("E30").Value = '  -2.40%  '
$ws.Range("D31").Value = '165.35'
$ws.Range("E31").Value = '  +1.28%  '
$ws.Range("D32").Value = '6.12'
$ws.Range("E32").Value = '  +4.81%  '
$ws.Range("D33").Value = '0.0885'
$ws.Range("E33").Value = '  +0.77%  '
$ws.Range("E34").Value = '  +0.88%  '
$ws.Range("E35").Value = '  -0.91%  '
$ws.Range("E36").Value = '  -13.25%  '
$ws.Range("D37").Value = '4.62'
$ws.Range("E37").Value = '  +0.84%  '
$ws.Range("D38").Value = '0.0355'
$ws.Range("E38").Value = '  +1.98%  '
$ws.Range("D39").Value = '3.76'
$ws.Range("E39").Value = '  +4.52%  '
$ws.Range("E40").Value = '  -5.29%  '
$ws.Range("D41").Value = '1.54'
$ws.Range("E41").Value = '  +5.39%  '
$ws.Range("D42").Value = '70.21'
$ws.Range("E42").Value = '  +0.30%  '
$ws.Range("D43").Value = '94.94'
$ws.Range("E43").Value = '  -6.84%  '
$ws.Range("E44").Value = '  -0.16%  '
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("D46").Value = '12.44'
$ws.Range("E46").Value = '  +4.15%  '
$ws.Range("D47").Value = '81.11'
$ws.Range("E47").Value = '  +5.90%  '
$ws.Range("D48").Value = '113.46'
$ws.Range("E48").Value = '  -1.99%  '
$ws.Range("D49").Value = '8.90'
$ws.Range("E49").Value = '  -1.74%  '
$ws.Range("D50").Value = '5.29'
$ws.Range("E50").Value = '  -0.44%  '
$ws.Range("D51").Value = '1.584.15'
$ws.Range("E51").Value = '  +2.28%  '
